# Apply updated cryptocurrency price/volume data to columns D (Price) and E (Volume(1h)).
# Both columns hold plain text in the source sheet (inline strings), so values that
# look like plain numbers ("1.01", "301.68", ...) are written with a leading apostrophe
# to force Excel to keep them as text, then the cell style is reset to match a sibling
# text cell so no stray number-format / quote-prefix styling is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$textStyleRef = $ws.Range("D2")

$ws.Range("D2").Value = "43.816.74"
$ws.Range("E2").Value = "  -0.39%  "

$ws.Range("D3").Value = "2.221.91"
$ws.Range("E3").Value = "  -1.53%  "

$ws.Range("D4").Value = "'1.01"
$ws.Range("D4").Style = $textStyleRef.Style
$ws.Range("E4").Value = "  +0.40%  "

$ws.Range("D5").Value = "'301.68"
$ws.Range("D5").Style = $textStyleRef.Style
$ws.Range("E5").Value = "  -5.42%  "

$ws.Range("D6").Value = "'92.46"
$ws.Range("D6").Style = $textStyleRef.Style
$ws.Range("E6").Value = "  -9.25%  "

$ws.Range("D7").Value = "'0.562"
$ws.Range("D7").Style = $textStyleRef.Style
$ws.Range("E7").Value = "  -2.67%  "

$ws.Range("E8").Value = "  +0.24%  "

$ws.Range("D9").Value = "'0.511"
$ws.Range("D9").Style = $textStyleRef.Style
$ws.Range("E9").Value = "  -8.32%  "

$ws.Range("D10").Value = "'33.63"
$ws.Range("D10").Style = $textStyleRef.Style
$ws.Range("E10").Value = "  -9.77%  "

$ws.Range("D11").Value = "'0.0790"
$ws.Range("D11").Style = $textStyleRef.Style
$ws.Range("E11").Value = "  -5.16%  "

$ws.Range("D12").Value = "'7.00"
$ws.Range("D12").Style = $textStyleRef.Style
$ws.Range("E12").Value = "  -9.05%  "

$ws.Range("E13").Value = "  -3.51%  "

$ws.Range("D14").Value = "2.560.79"
$ws.Range("E14").Value = "  -1.28%  "

$ws.Range("D15").Value = "2.251.50"
$ws.Range("E15").Value = "  -0.07%  "

$ws.Range("D16").Value = "'0.802"
$ws.Range("D16").Style = $textStyleRef.Style
$ws.Range("E16").Value = "  -6.87%  "

$ws.Range("D17").Value = "'13.31"
$ws.Range("D17").Style = $textStyleRef.Style
$ws.Range("E17").Value = "  -6.65%  "

$ws.Range("D18").Value = "43.594.77"
$ws.Range("E18").Value = "  -0.64%  "

$ws.Range("D19").Value = "0.0₃0941"
$ws.Range("E19").Value = "  -4.74%  "

$ws.Range("D20").Value = "'11.96"
$ws.Range("D20").Style = $textStyleRef.Style
$ws.Range("E20").Value = "  -11.27%  "

$ws.Range("D21").Value = "'6.06"
$ws.Range("D21").Style = $textStyleRef.Style
$ws.Range("E21").Value = "  -7.54%  "

$ws.Range("D22").Value = "'63.68"
$ws.Range("D22").Style = $textStyleRef.Style
$ws.Range("E22").Value = "  -3.22%  "

$ws.Range("D23").Value = "'233.03"
$ws.Range("D23").Style = $textStyleRef.Style
$ws.Range("E23").Value = "  -1.28%  "

$ws.Range("E24").Value = "  -8.61%  "

$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("D26").Value = "'1.89"
$ws.Range("D26").Style = $textStyleRef.Style
$ws.Range("E26").Value = "  -11.27%  "

$ws.Range("D27").Value = "'9.60"
$ws.Range("D27").Style = $textStyleRef.Style
$ws.Range("E27").Value = "  -5.68%  "

$ws.Range("E28").Value = "  -2.34%  "

$ws.Range("D29").Value = "'35.28"
$ws.Range("D29").Style = $textStyleRef.Style
$ws.Range("E29").Value = "  -5.89%  "

$ws.Range("D30").Value = "'5.77"
$ws.Range("D30").Style = $textStyleRef.Style
$ws.Range("E30").Value = "  -7.95%  "

$ws.Range("D31").Value = "'19.52"
$ws.Range("D31").Style = $textStyleRef.Style
$ws.Range("E31").Value = "  -3.60%  "

$ws.Range("D32").Value = "'149.93"
$ws.Range("D32").Style = $textStyleRef.Style
$ws.Range("E32").Value = "  -5.11%  "

$ws.Range("D33").Value = "'0.0788"
$ws.Range("D33").Style = $textStyleRef.Style
$ws.Range("E33").Value = "  -7.83%  "

$ws.Range("E34").Value = "  -4.19%  "

$ws.Range("D35").Value = "'3.17"
$ws.Range("D35").Style = $textStyleRef.Style
$ws.Range("E35").Value = "  +2.20%  "

$ws.Range("D36").Value = "'0.116"
$ws.Range("D36").Style = $textStyleRef.Style
$ws.Range("E36").Value = "  -2.21%  "

$ws.Range("D37").Value = "'0.104"
$ws.Range("D37").Style = $textStyleRef.Style
$ws.Range("E37").Value = "  -10.90%  "

$ws.Range("D38").Value = "'1.72"
$ws.Range("D38").Style = $textStyleRef.Style
$ws.Range("E38").Value = "  -12.63%  "

$ws.Range("D39").Value = "'14.25"
$ws.Range("D39").Style = $textStyleRef.Style
$ws.Range("E39").Value = "  -11.82%  "

$ws.Range("D40").Value = "'3.70"
$ws.Range("D40").Style = $textStyleRef.Style
$ws.Range("E40").Value = "  -12.37%  "

$ws.Range("D41").Value = "'0.0293"
$ws.Range("D41").Style = $textStyleRef.Style
$ws.Range("E41").Value = "  -7.68%  "

$ws.Range("D42").Value = "'3.20"
$ws.Range("D42").Style = $textStyleRef.Style
$ws.Range("E42").Value = "  -14.63%  "

$ws.Range("E43").Value = "  +0.38%  "

$ws.Range("D44").Value = "1.719.90"
$ws.Range("E44").Value = "  -4.54%  "

$ws.Range("D45").Value = "'81.91"
$ws.Range("D45").Style = $textStyleRef.Style
$ws.Range("E45").Value = "  -1.14%  "

$ws.Range("D46").Value = "'4.87"
$ws.Range("D46").Style = $textStyleRef.Style
$ws.Range("E46").Value = "  -6.85%  "

$ws.Range("D47").Value = "'0.182"
$ws.Range("D47").Style = $textStyleRef.Style
$ws.Range("E47").Value = "  -8.81%  "

$ws.Range("D48").Value = "'97.90"
$ws.Range("D48").Style = $textStyleRef.Style
$ws.Range("E48").Value = "  -5.86%  "

$ws.Range("D49").Value = "'7.97"
$ws.Range("D49").Style = $textStyleRef.Style
$ws.Range("E49").Value = "  -5.23%  "

$ws.Range("D50").Value = "'52.83"
$ws.Range("D50").Style = $textStyleRef.Style
$ws.Range("E50").Value = "  -10.16%  "

$ws.Range("D51").Value = "'66.09"
$ws.Range("D51").Style = $textStyleRef.Style
$ws.Range("E51").Value = "  -13.02%  "
